$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Z2:Z45").Value = "2025-10-17T07:09:32.480043"
$ws.Range("Z46:Z74").Value = "2025-10-17T07:09:32.592752"
$ws.Range("Z75").Value = "2025-10-17T07:09:32.713135"
$ws.Range("Z76:Z83").Value = "2025-10-17T07:09:32.714135"
$ws.Range("Z84:Z94").Value = "2025-10-17T07:09:32.715135"
$ws.Range("Z95:Z102").Value = "2025-10-17T07:09:32.716135"
$ws.Range("Z103:Z112").Value = "2025-10-17T07:09:32.830037"
